$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2: update D2, E2
Set-TextValue $ws.Range("D2") '28.222.19'
Set-TextValue $ws.Range("E2") '  +1.74%  '

# Row 3: update D3, E3
Set-TextValue $ws.Range("D3") '1.802.40'
Set-TextValue $ws.Range("E3") '  +2.48%  '

# Row 4: update D4, E4
Set-TextValue $ws.Range("D4") '1.001'
Set-TextValue $ws.Range("E4") '  -0.34%  '

# Row 5: update D5, E5
Set-TextValue $ws.Range("D5") '324.64'
Set-TextValue $ws.Range("E5") '  -0.69%  '

# Row 6: update D6, E6
Set-TextValue $ws.Range("D6") '1.000'
Set-TextValue $ws.Range("E6") '  -0.07%  '

# Row 7: update D7, E7
Set-TextValue $ws.Range("D7") '0.4301'
Set-TextValue $ws.Range("E7") '  -3.04%  '

# Row 8: update D8, E8
Set-TextValue $ws.Range("D8") '0.3638'
Set-TextValue $ws.Range("E8") '  -3.07%  '

# Row 9: update D9, E9
Set-TextValue $ws.Range("D9") '44.73'
Set-TextValue $ws.Range("E9") '  -1.66%  '

# Row 10: update D10, E10
Set-TextValue $ws.Range("D10") '0.07575'
Set-TextValue $ws.Range("E10") '  -0.98%  '

# Row 11: update D11, E11
Set-TextValue $ws.Range("D11") '1.126'
Set-TextValue $ws.Range("E11") '  +0.04%  '

# Row 12: update D12, E12
Set-TextValue $ws.Range("D12") '1.000'
Set-TextValue $ws.Range("E12") '  -0.05%  '

# Row 13: update D13, E13
Set-TextValue $ws.Range("D13") '21.78'
Set-TextValue $ws.Range("E13") '  +0.21%  '

# Row 14: update D14, E14
Set-TextValue $ws.Range("D14") '6.209'
Set-TextValue $ws.Range("E14") '  +0.05%  '

# Row 15: update D15, E15
Set-TextValue $ws.Range("D15") '7.390'
Set-TextValue $ws.Range("E15") '  -0.87%  '

# Row 16: update D16, E16
Set-TextValue $ws.Range("D16") '1.814.59'
Set-TextValue $ws.Range("E16") '  +3.29%  '

# Row 17: update D17, E17
Set-TextValue $ws.Range("D17") '93.16'
Set-TextValue $ws.Range("E17") '  +4.68%  '

# Row 18: update D18, E18
Set-TextValue $ws.Range("D18") '0.00001071'
Set-TextValue $ws.Range("E18") '  -0.50%  '

# Row 19: update D19, E19
Set-TextValue $ws.Range("D19") '0.06379'
Set-TextValue $ws.Range("E19") '  +2.67%  '

# Row 20: update D20, E20
Set-TextValue $ws.Range("D20") '1.001'
Set-TextValue $ws.Range("E20") '  +0.00%  '

# Row 21: update D21
Set-TextValue $ws.Range("D21") '17.30'

# Row 22: update D22, E22
Set-TextValue $ws.Range("D22") '6.015'
Set-TextValue $ws.Range("E22") '  -2.71%  '

# Row 23: update D23, E23
Set-TextValue $ws.Range("D23") '28.223.09'
Set-TextValue $ws.Range("E23") '  +1.64%  '

# Row 24: update D24, E24
Set-TextValue $ws.Range("D24") '11.44'
Set-TextValue $ws.Range("E24") '  -1.67%  '

# Row 25: update D25, E25
Set-TextValue $ws.Range("D25") '2.173'
Set-TextValue $ws.Range("E25") '  -6.15%  '

# Row 26: update D26, E26
Set-TextValue $ws.Range("D26") '160.14'
Set-TextValue $ws.Range("E26") '  +4.06%  '

# Row 27: update D27, E27
Set-TextValue $ws.Range("D27") '20.45'
Set-TextValue $ws.Range("E27") '  -1.40%  '

# Row 28: update D28, E28
Set-TextValue $ws.Range("D28") '2.018.62'
Set-TextValue $ws.Range("E28") '  +3.23%  '

# Row 29: update D29, E29
Set-TextValue $ws.Range("D29") '2.245'
Set-TextValue $ws.Range("E29") '  -5.03%  '

# Row 30: update D30, E30
Set-TextValue $ws.Range("D30") '128.37'
Set-TextValue $ws.Range("E30") '  +0.03%  '

# Row 31: update D31, E31
Set-TextValue $ws.Range("D31") '1.182'
Set-TextValue $ws.Range("E31") '  -2.96%  '

# Row 32: update D32, E32
Set-TextValue $ws.Range("D32") '5.903'
Set-TextValue $ws.Range("E32") '  +2.30%  '

# Row 33: update D33, E33
Set-TextValue $ws.Range("D33") '0.09047'
Set-TextValue $ws.Range("E33") '  -3.40%  '

# Row 34: update D34, E34
Set-TextValue $ws.Range("D34") '3.530'
Set-TextValue $ws.Range("E34") '  -3.56%  '

# Row 35: update D35, E35
Set-TextValue $ws.Range("D35") '12.85'
Set-TextValue $ws.Range("E35") '  +0.94%  '

# Row 36: update D36, E36
Set-TextValue $ws.Range("D36") '0.02365'
Set-TextValue $ws.Range("E36") '  +1.70%  '

# Row 37: update D37
Set-TextValue $ws.Range("D37") '5.150'

# Row 38: update D38, E38
Set-TextValue $ws.Range("D38") '0.6518'
Set-TextValue $ws.Range("E38") '  +0.29%  '

# Row 39: update B39, C39, D39, E39
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D39") '0.2131'
Set-TextValue $ws.Range("E39") '  -2.56%  '

# Row 40: update B40, C40, D40, E40
$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D40") '0.06135'
Set-TextValue $ws.Range("E40") '  -0.34%  '

# Row 41: update D41, E41
Set-TextValue $ws.Range("D41") '1.198'
Set-TextValue $ws.Range("E41") '  -0.52%  '

# Row 42: update E42
Set-TextValue $ws.Range("E42") '  +0.95%  '

# Row 43: update D43, E43
Set-TextValue $ws.Range("D43") '7.978'
Set-TextValue $ws.Range("E43") '  -0.41%  '

# Row 44: update D44, E44
Set-TextValue $ws.Range("D44") '1.000'
Set-TextValue $ws.Range("E44") '  -0.02%  '

# Row 45: update D45, E45
Set-TextValue $ws.Range("D45") '13.60'
Set-TextValue $ws.Range("E45") '  -1.37%  '

# Row 46: update D46, E46
Set-TextValue $ws.Range("D46") '0.6038'
Set-TextValue $ws.Range("E46") '  +0.28%  '

# Row 47: update D47, E47
Set-TextValue $ws.Range("D47") '3.712'
Set-TextValue $ws.Range("E47") '  -1.46%  '

# Row 48: update D48, E48
Set-TextValue $ws.Range("D48") '125.71'
Set-TextValue $ws.Range("E48") '  -0.56%  '

# Row 49: update D49, E49
Set-TextValue $ws.Range("D49") '1.995'
Set-TextValue $ws.Range("E49") '  -0.23%  '

# Row 50: update D50, E50
Set-TextValue $ws.Range("D50") '1.161'
Set-TextValue $ws.Range("E50") '  +2.05%  '

# Row 51: update D51
Set-TextValue $ws.Range("D51") '0.06975'
